# Apply the edits described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the extraction timestamp in the study description.
$ws.Range("A2").Value = "This is an interesting study.Extracted on : 2022/09/26 11:58:55"

# Rename variable label.
$ws.Range("A7").Value = "Compactness quartile"

# Update OR (univariate) column (E) values.
$ws.Range("E5").Value = "2.811 (2.372-3.42, p=1e-28)"
$ws.Range("E6").Value = "1.264 (1.203-1.333, p=3e-19)"
$ws.Range("E8").Value = "2.709 (1.237-6.416, p=0.02)"
$ws.Range("E9").Value = "12.475 (6.162-28.139, p=5e-11)"
$ws.Range("E10").Value = "65.932 (31.127-155.485, p=8e-25)"

# Update OR (model 1) column (F) values.
$ws.Range("F5").Value = "4.043 (3.073-5.638, p=1e-19)"
$ws.Range("F6").Value = "1.331 (1.258-1.414, p=7e-22)"
$ws.Range("F8").Value = "3.4 (1.396-8.974, p=0.009)"
$ws.Range("F9").Value = "12.686 (5.422-32.866, p=3e-08)"
$ws.Range("F10").Value = "54.214 (20.658-157.555, p=1e-14)"

# Update OR (model 2) column (G) values.
$ws.Range("G5").Value = "2.716 (2.204-3.46, p=3e-18)"
$ws.Range("G6").Value = "1.246 (1.162-1.341, p=1e-09)"
$ws.Range("G8").Value = "0.767 (0.314-1.964, p=0.6)"
$ws.Range("G9").Value = "0.524 (0.184-1.516, p=0.2)"
$ws.Range("G10").Value = "0.219 (0.051-0.908, p=0.04)"
